$wb = $excel.ActiveWorkbook

# --- Rename Sheet3 -> 20120309 and populate it ---
$s3 = $wb.Worksheets.Item(3)
$s3.Name = "20120309"

$s3.Cells.Item(10, 6).Value = 1
$s3.Cells.Item(10, 7).Value = "SOA study"
$s3.Cells.Item(10, 9).Value = "16:00-16:30"
$s3.Cells.Item(10, 11).Value = "0.5h"

$s3.Cells.Item(11, 6).Value = 2
$s3.Cells.Item(11, 7).Value = "Hudson study"
$s3.Cells.Item(11, 9).Value = "16:30-18:00"
$s3.Cells.Item(11, 11).Value = "1.5h"

$s3.Columns.Item(7).ColumnWidth = 11.660714285714286

# --- Add 20120725 ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$s4 = $wb.Worksheets.Add($null, $last)
$s4.Name = "20120725"

$s4.Cells.Item(13, 7).Value = 1
$s4.Cells.Item(13, 8).Value = "firebird study"
$s4.Cells.Item(13, 11).Value = "15:00-17:30"
$s4.Cells.Item(13, 13).Value = "2.5h"

$s4.PageSetup.PaperSize = 9
$s4.PageSetup.Orientation = 1

# --- Add 20120727 ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$s5 = $wb.Worksheets.Add($null, $last)
$s5.Name = "20120727"

$s5.Cells.Item(13, 7).Value = 1
$s5.Cells.Item(13, 8).Value = "firebird study"
$s5.Cells.Item(13, 11).Value = "10:00-12:00"
$s5.Cells.Item(13, 13).Value = "2h"

# --- Add 20120728 ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$s6 = $wb.Worksheets.Add($null, $last)
$s6.Name = "20120728"

$s6.Cells.Item(11, 5).Value = 1
$s6.Cells.Item(11, 6).Value = "change BaseAction  to DispatchAction"
$s6.Cells.Item(11, 11).Value = "11:00-12:30"
$s6.Cells.Item(11, 13).Value = "1.5h"

$s6.Columns.Item(10).ColumnWidth = 24.410714285714285

$s6.PageSetup.PaperSize = 9
$s6.PageSetup.Orientation = 1

# --- Selections / active cells per final sheet views ---
$s2 = $wb.Worksheets.Item(2)
$s2.Activate()
[void]$s2.Range("K7").Select()

[void]$s3.Activate()
[void]$s3.Range("H17").Select()

[void]$s4.Activate()
[void]$s4.Range("G13:M13").Select()

[void]$s5.Activate()
[void]$s5.Range("I15").Select()

[void]$s6.Activate()
[void]$s6.Range("F11").Select()
